$d = $word.ActiveDocument

$replacements = @(
    @("N = 1,838", "N = 1,790"),
    @("N = 1,841", "N = 1,793"),
    @("N = 1,848", "N = 1,800"),
    @("N = 1,874", "N = 1,825"),
    @("4.82 (0.23)", "4.82 (0.22)"),
    @("10.90 (0.33)", "10.91 (0.33)"),
    @("937 / 1,838 (51%)", "917 / 1,790 (51%)"),
    @("939 / 1,841 (51%)", "919 / 1,793 (51%)"),
    @("945 / 1,848 (51%)", "925 / 1,800 (51%)"),
    @("955 / 1,874 (51%)", "935 / 1,825 (51%)"),
    @("901 / 1,838 (49%)", "873 / 1,790 (49%)"),
    @("902 / 1,841 (49%)", "874 / 1,793 (49%)"),
    @("903 / 1,848 (49%)", "875 / 1,800 (49%)"),
    @("919 / 1,874 (49%)", "890 / 1,825 (49%)"),
    @("35 / 1,838 (1.9%)", "32 / 1,790 (1.8%)"),
    @("35 / 1,841 (1.9%)", "32 / 1,793 (1.8%)"),
    @("35 / 1,848 (1.9%)", "32 / 1,800 (1.8%)"),
    @("37 / 1,874 (2.0%)", "33 / 1,825 (1.8%)"),
    @("0.26 (0.95)", "0.27 (0.95)"),
    @("0.23 (0.96)", "0.24 (0.96)"),
    @("0.21 (0.97)", "0.22 (0.97)"),
    @("932 (538)", "932 (541)"),
    @("1,038 (553)", "1,036 (550)"),
    @("1,250 (601)", "1,248 (601)"),
    @("1,375 (643)", "1,371 (644)"),
    @("62 (34)", "63 (34)"),
    @("553 (81)", "552 (81)"),
    @("1,632 / 1,833 (89%)", "1,591 / 1,787 (89%)"),
    @("201 / 1,833 (11%)", "196 / 1,787 (11%)"),
    @("330 / 1,745 (19%)", "322 / 1,704 (19%)"),
    @("657 / 1,745 (38%)", "639 / 1,704 (38%)"),
    @("502 / 1,745 (29%)", "492 / 1,704 (29%)"),
    @("216 / 1,745 (12%)", "212 / 1,704 (12%)"),
    @("40 / 1,745 (2.3%)", "39 / 1,704 (2.3%)"),
    @("0.31 (0.99)", "0.30 (0.99)"),
    @("4.06 (0.65)", "4.07 (0.65)"),
    @("1.38 (0.12)", "1.39 (0.12)"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $old
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = $new
    $result = $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $old"
    }
}

Write-Host "Done"
